$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8962752223014832
$ws.Range("B1").Value = 1.223539352416992
$ws.Range("C1").Value = 2.062266826629639
$ws.Range("D1").Value = 4.513389110565186
$ws.Range("E1").Value = 2.172746658325195
